$wb = $excel.ActiveWorkbook

# 1. Fix the typo in the first sheet's name: "TA-Cours-Groups" -> "TA-Course-Groups"
$wb.Sheets.Item(1).Name = "TA-Course-Groups"

# 2. On the "Course-Groups" sheet, wrap + grow row 3 (course name is long) and
#    move that sheet's selection to A3. A new cell style (wrap + left align)
#    gets created automatically by the engine when WrapText is applied.
$ws = $wb.Sheets.Item("Course-Groups")
$ws.Range("A3").WrapText = $true
$ws.Rows.Item(3).RowHeight = 37.5
$ws.Range("A3").Select()

# Restore the originally active sheet/tab (selecting A3 above would otherwise
# shift which tab is marked as selected).
$wb.Sheets.Item(1).Activate()
